$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 542; this shifts the existing rows 542:573
# down to 543:574 (values, formatting, etc. all move with the rows).
$ws.Rows.Item(542).Insert()

# Populate the newly inserted row 542 with the new weekly record.
$ws.Cells.Item(542, 1).Value = 10
$ws.Cells.Item(542, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(542, 3).Value = "La Araucanía"
$ws.Cells.Item(542, 4).Value = 45041
$ws.Cells.Item(542, 5).Value = 9
$ws.Cells.Item(542, 6).Value = 100114014
$ws.Cells.Item(542, 7).Value = "Betarraga"
$ws.Cells.Item(542, 8).Value = "Sin especificar"
$ws.Cells.Item(542, 9).Value = "Primera"
$ws.Cells.Item(542, 10).Value = 40
$ws.Cells.Item(542, 11).Value = 10000
$ws.Cells.Item(542, 12).Value = 10000
$ws.Cells.Item(542, 13).Value = 10000
$ws.Cells.Item(542, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(542, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(542, 16).Value = 833
$ws.Cells.Item(542, 17).Value = 12
$ws.Cells.Item(542, 18).Value = "Hortaliza"
